$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.738.04"
$ws.Range("E2").Value = "  +2.47%  "
$ws.Range("D3").Value = "'1.894.49"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D5").Value = "'246.71"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D7").Value = "'0.4932"
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("D8").Value = "'0.2959"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "'0.06816"
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("D10").Value = "'1.894.66"
$ws.Range("D11").Value = "'17.32"
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("D12").Value = "'92.25"
$ws.Range("E12").Value = "  +6.86%  "
$ws.Range("D13").Value = "'0.07262"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "'0.6841"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "'5.101"
$ws.Range("E15").Value = "  +4.60%  "
$ws.Range("D16").Value = "'30.713.60"
$ws.Range("D17").Value = "'0.000007985"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("E18").Value = "  +4.47%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'2.140.09"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "'4.856"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").Value = "'190.11"
$ws.Range("E23").Value = "  +34.56%  "
$ws.Range("D24").Value = "'6.078"
$ws.Range("E24").Value = "  +7.62%  "
$ws.Range("D25").Value = "'9.409"
$ws.Range("E25").Value = "  +3.86%  "
$ws.Range("D26").Value = "'155.45"
$ws.Range("E26").Value = "  +4.22%  "
$ws.Range("E27").Value = "  +12.23%  "
$ws.Range("D28").Value = "'1.931"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "'1.401"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").Value = "'4.394"
$ws.Range("E30").Value = "  +5.27%  "
$ws.Range("E31").Value = "  +2.72%  "
$ws.Range("D32").Value = "'4.057"
$ws.Range("E32").Value = "  +2.86%  "
$ws.Range("D33").Value = "'0.05207"
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("D34").Value = "'0.7494"
$ws.Range("D35").Value = "'1.129"
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("D36").Value = "'2.713"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").Value = "'0.01870"
$ws.Range("D38").Value = "'2.677"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Value = "'2.169"
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "'0.9384"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  +4.65%  "
$ws.Range("E42").Value = "  +4.21%  "
$ws.Range("D43").Value = "'5.825"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "'7.723"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("D47").Value = "'0.05859"
$ws.Range("E47").Value = "  +3.70%  "
$ws.Range("D48").Value = "'8.821"
$ws.Range("E48").Value = "  +7.30%  "
$ws.Range("D49").Value = "'0.3974"
$ws.Range("E49").Value = "  +5.81%  "
$ws.Range("D50").Value = "'1.426"
$ws.Range("E50").Value = "  +6.88%  "
$ws.Range("D51").Value = "'33.63"

Write-Host "cryptos list updated"
